# "Generate Report for Handoff"
# Moves the localization status from "In Translation" to "Ready for handoff"
# and refreshes the associated handoff timestamps, on all three sheets
# (Overview, zh-cn, de-de). Also widens the Status-related columns that
# grew to fit the new, longer status text.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-19 00:56:03"
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-19 00:55:56"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-19 00:56:03"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
